$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 224
$ws.Range("F6").Value = 77
$ws.Range("F7").Value = 794
$ws.Range("F8").Value = 453
$ws.Range("F9").Value = 62
$ws.Range("F10").Value = 277
$ws.Range("F14").Value = 395
$ws.Range("F15").Value = 6393
$ws.Range("F18").Value = 16
$ws.Range("F19").Value = 7354
$ws.Range("F22").Value = 3330
$ws.Range("F23").Value = 762
$ws.Range("F24").Value = 843
$ws.Range("F25").Value = 4492
$ws.Range("F26").Value = 338
$ws.Range("F27").Value = 175
$ws.Range("F28").Value = 169
$ws.Range("F29").Value = 1386
$ws.Range("F30").Value = 133
$ws.Range("F31").Value = 44
$ws.Range("F32").Value = 11
$ws.Range("F33").Value = 1079
$ws.Range("F34").Value = 1494
$ws.Range("F35").Value = 2105

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 17

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 235
$ws.Range("F3").Value = 1191

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 235
$ws.Range("F4").Value = 1191
$ws.Range("F8").Value = 224
$ws.Range("F9").Value = 77
$ws.Range("F10").Value = 794
$ws.Range("F11").Value = 453
$ws.Range("F12").Value = 62
$ws.Range("F13").Value = 277
$ws.Range("F18").Value = 395
$ws.Range("F19").Value = 6393
$ws.Range("F22").Value = 16
$ws.Range("F23").Value = 7354
$ws.Range("F26").Value = 3330
$ws.Range("F27").Value = 762
$ws.Range("F28").Value = 843
$ws.Range("F29").Value = 4492
$ws.Range("F30").Value = 338
$ws.Range("F32").Value = 175
$ws.Range("F33").Value = 169
$ws.Range("F34").Value = 1386
$ws.Range("F35").Value = 133
$ws.Range("F36").Value = 44
$ws.Range("F37").Value = 11
$ws.Range("F38").Value = 1079
$ws.Range("F39").Value = 1494
$ws.Range("F40").Value = 17
$ws.Range("F41").Value = 2105
